$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    "53.06.13.2021",
    "53.06.13.2020",
    "53.06.13.2019",
    "53.06.13.2018",
    "53.06.13.2017",
    "53.06.13.2016",
    "53.06.13.2015",
    "53.06.13.2014"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("G9").Select()
